$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ciudad" column (old column F). This shifts the columns that
# followed it (destinatario_telefono, destinatario_correo, peso, ...) one
# position to the left, taking their values/styles/widths with them.
$ws.Columns("F").Delete()

# Insert a brand new column before the (now shifted) "peso" column so that a
# new "unidades" column appears at H, while "peso" moves back out to I.
$ws.Columns("H").Insert()

# Populate the new "unidades" column.
$ws.Range("H1").Value = "unidades"
$ws.Range("H2").Value = 10

# Give the new column roughly the same width as its neighbours (~12.1 chars).
$ws.Columns("H").ColumnWidth = 11.26

# Update the active selection to H2, as recorded in the saved view state.
$ws.Range("H2").Select()

# Update header/footer font style from "Normal" to "Regular".
$ws.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12&A'
$ws.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12Página &P'
